$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("X2").Value = "Utility (Percent)"
$ws.Range("L3").Value = "55 msec"
$ws.Range("M3").Value = "2542 msec"
$ws.Range("N3").Value = "1266.26 msec"
$ws.Range("O3").Value = "1334 usec"
$ws.Range("P3").Value = "2330.9k usec"
$ws.Range("Q3").Value = "30491.43 usec"
$ws.Range("L4").Value = "587 usec"
$ws.Range("M4").Value = "5288.5k usec"
$ws.Range("N4").Value = "920981.53 usec"
$ws.Range("O4").Value = "2 msec"
$ws.Range("P4").Value = "7324 msec"
$ws.Range("Q4").Value = "1140.99 msec"
$ws.Range("L5").Value = "61 msec"
$ws.Range("M5").Value = "2436 msec"
$ws.Range("N5").Value = "1271.57 msec"
$ws.Range("O5").Value = "1980 usec"
$ws.Range("P5").Value = "1609.6k usec"
$ws.Range("Q5").Value = "17412.80 usec"
$ws.Range("L6").Value = "70 msec"
$ws.Range("M6").Value = "2949 msec"
$ws.Range("N6").Value = "1251.44 msec"
$ws.Range("O6").Value = "1278 usec"
$ws.Range("P6").Value = "2624.3k usec"
$ws.Range("Q6").Value = "65114.86 usec"
$ws.Range("L7").Value = "375 usec"
$ws.Range("M7").Value = "5614.5k usec"
$ws.Range("N7").Value = "674565.98 usec"
$ws.Range("O7").Value = "2 msec"
$ws.Range("P7").Value = "8941 msec"
$ws.Range("Q7").Value = "2003.07 msec"
$ws.Range("L8").Value = "645 usec"
$ws.Range("M8").Value = "5623.4k usec"
$ws.Range("N8").Value = "716647.92 usec"
$ws.Range("O8").Value = "6 msec"
$ws.Range("P8").Value = "6499 msec"
$ws.Range("Q8").Value = "1892.25 msec"
$ws.Range("L9").Value = "769 usec"
$ws.Range("M9").Value = "12654k usec"
$ws.Range("N9").Value = "1438889.67 usec"
$ws.Range("O9").Value = "2 msec"
$ws.Range("P9").Value = "14574 msec"
$ws.Range("Q9").Value = "1386.30 msec"
$ws.Range("L10").Value = "37 msec"
$ws.Range("M10").Value = "3340 msec"
$ws.Range("N10").Value = "1264.56 msec"
$ws.Range("O10").Value = "1508 usec"
$ws.Range("P10").Value = "3126.8k usec"
$ws.Range("Q10").Value = "34680.21 usec"
$ws.Range("L11").Value = "17 msec"
$ws.Range("M11").Value = "2895 msec"
$ws.Range("N11").Value = "1276.04 msec"
$ws.Range("O11").Value = "1924 usec"
$ws.Range("P11").Value = "1204.6k usec"
$ws.Range("Q11").Value = "7259.67 usec"
$ws.Range("L12").Value = "1934 usec"
$ws.Range("M12").Value = "7421.7k usec"
$ws.Range("N12").Value = "1230073.08 usec"
$ws.Range("O12").Value = "3 msec"
$ws.Range("P12").Value = "7280 msec"
$ws.Range("Q12").Value = "142.84 msec"
$ws.Range("L13").Value = "6 msec"
$ws.Range("M13").Value = "3987 msec"
$ws.Range("N13").Value = "1213.73 msec"
$ws.Range("O13").Value = "2 msec"
$ws.Range("P13").Value = "4351 msec"
$ws.Range("Q13").Value = "155.15 msec"
$ws.Range("L14").Value = "14 msec"
$ws.Range("M14").Value = "3613 msec"
$ws.Range("N14").Value = "1266.64 msec"
$ws.Range("O14").Value = "3 msec"
$ws.Range("P14").Value = "2811 msec"
$ws.Range("Q14").Value = "28.67 msec"
$ws.Range("L15").Value = "612 usec"
$ws.Range("M15").Value = "5503.2k usec"
$ws.Range("N15").Value = "688017.37 usec"
$ws.Range("O15").Value = "4 msec"
$ws.Range("P15").Value = "10922 msec"
$ws.Range("Q15").Value = "1972.19 msec"
$ws.Range("L16").Value = "53 msec"
$ws.Range("M16").Value = "2644 msec"
$ws.Range("N16").Value = "1274.94 msec"
$ws.Range("O16").Value = "1480 usec"
$ws.Range("P16").Value = "1420.8k usec"
$ws.Range("Q16").Value = "9977.52 usec"
$ws.Range("L17").Value = "885 usec"
$ws.Range("M17").Value = "3884.8k usec"
$ws.Range("N17").Value = "1179691.54 usec"
$ws.Range("O17").Value = "1676 usec"
$ws.Range("P17").Value = "4202.3k usec"
$ws.Range("Q17").Value = "236328.47 usec"
$ws.Range("L18").Value = "72 msec"
$ws.Range("M18").Value = "2907 msec"
$ws.Range("N18").Value = "1244.99 msec"
$ws.Range("O18").Value = "1631 usec"
$ws.Range("P18").Value = "2015.4k usec"
$ws.Range("Q18").Value = "80412.84 usec"
$ws.Range("L19").Value = "55 msec"
$ws.Range("M19").Value = "2435 msec"
$ws.Range("N19").Value = "1270.10 msec"
$ws.Range("O19").Value = "1738 usec"
$ws.Range("P19").Value = "3516.5k usec"
$ws.Range("Q19").Value = "20959.92 usec"
$ws.Range("L20").Value = "11 msec"
$ws.Range("M20").Value = "2733 msec"
$ws.Range("N20").Value = "1270.19 msec"
$ws.Range("O20").Value = "3 msec"
$ws.Range("P20").Value = "1750 msec"
$ws.Range("Q20").Value = "20.55 msec"
$ws.Range("L21").Value = "477 usec"
$ws.Range("M21").Value = "5528.2k usec"
$ws.Range("N21").Value = "1012049.52 usec"
$ws.Range("O21").Value = "1869 usec"
$ws.Range("P21").Value = "7415.9k usec"
$ws.Range("Q21").Value = "884699.58 usec"
$ws.Range("L22").Value = "74 msec"
$ws.Range("M22").Value = "2527 msec"
$ws.Range("N22").Value = "1274.43 msec"
$ws.Range("O22").Value = "2 msec"
$ws.Range("P22").Value = "1636 msec"
$ws.Range("Q22").Value = "11.09 msec"
$ws.Range("L23").Value = "52 msec"
$ws.Range("M23").Value = "2445 msec"
$ws.Range("N23").Value = "1236.71 msec"
$ws.Range("O23").Value = "1636 usec"
$ws.Range("P23").Value = "2706.6k usec"
$ws.Range("Q23").Value = "99251.99 usec"
